$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.046640396118164
$ws.Range("B1").Value = 6.173583984375
$ws.Range("C1").Value = 3.201253414154053
$ws.Range("D1").Value = 1.410980701446533
$ws.Range("E1").Value = 0.9910504817962646
